$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap rows 14 and 16 (columns A,B,D,E,F,G,H,Q,R) ---
$row14 = @{
    A = $ws.Range("A14").Value2
    B = $ws.Range("B14").Value2
    D = $ws.Range("D14").Value2
    E = $ws.Range("E14").Value2
    F = $ws.Range("F14").Value2
    G = $ws.Range("G14").Value2
    H = $ws.Range("H14").Value2
    Q = $ws.Range("Q14").Value2
    R = $ws.Range("R14").Value2
}

$row16 = @{
    A = $ws.Range("A16").Value2
    B = $ws.Range("B16").Value2
    D = $ws.Range("D16").Value2
    E = $ws.Range("E16").Value2
    F = $ws.Range("F16").Value2
    G = $ws.Range("G16").Value2
    H = $ws.Range("H16").Value2
    Q = $ws.Range("Q16").Value2
    R = $ws.Range("R16").Value2
}

$ws.Range("A14").Value = $row16.A
$ws.Range("B14").Value = $row16.B
$ws.Range("D14").Value = $row16.D
$ws.Range("E14").Value = $row16.E
$ws.Range("F14").Value = $row16.F
$ws.Range("G14").Value = $row16.G
$ws.Range("H14").Value = $row16.H
$ws.Range("Q14").Value = $row16.Q
$ws.Range("R14").Value = $row16.R

$ws.Range("A16").Value = $row14.A
$ws.Range("B16").Value = $row14.B
$ws.Range("D16").Value = $row14.D
$ws.Range("E16").Value = $row14.E
$ws.Range("F16").Value = $row14.F
$ws.Range("G16").Value = $row14.G
$ws.Range("H16").Value = $row14.H
$ws.Range("Q16").Value = $row14.Q
$ws.Range("R16").Value = $row14.R

# --- Rotate rows 19, 20, 21 (columns A,B,D,E,F,G,H,Q,R) ---
# new19 = old20, new20 = old21, new21 = old19
$row19 = @{
    A = $ws.Range("A19").Value2
    B = $ws.Range("B19").Value2
    D = $ws.Range("D19").Value2
    E = $ws.Range("E19").Value2
    F = $ws.Range("F19").Value2
    G = $ws.Range("G19").Value2
    H = $ws.Range("H19").Value2
    Q = $ws.Range("Q19").Value2
    R = $ws.Range("R19").Value2
}

$row20 = @{
    A = $ws.Range("A20").Value2
    B = $ws.Range("B20").Value2
    D = $ws.Range("D20").Value2
    E = $ws.Range("E20").Value2
    F = $ws.Range("F20").Value2
    G = $ws.Range("G20").Value2
    H = $ws.Range("H20").Value2
    Q = $ws.Range("Q20").Value2
    R = $ws.Range("R20").Value2
}

$row21 = @{
    A = $ws.Range("A21").Value2
    B = $ws.Range("B21").Value2
    D = $ws.Range("D21").Value2
    E = $ws.Range("E21").Value2
    F = $ws.Range("F21").Value2
    G = $ws.Range("G21").Value2
    H = $ws.Range("H21").Value2
    Q = $ws.Range("Q21").Value2
    R = $ws.Range("R21").Value2
}

$ws.Range("A19").Value = $row20.A
$ws.Range("B19").Value = $row20.B
$ws.Range("D19").Value = $row20.D
$ws.Range("E19").Value = $row20.E
$ws.Range("F19").Value = $row20.F
$ws.Range("G19").Value = $row20.G
$ws.Range("H19").Value = $row20.H
$ws.Range("Q19").Value = $row20.Q
$ws.Range("R19").Value = $row20.R

$ws.Range("A20").Value = $row21.A
$ws.Range("B20").Value = $row21.B
$ws.Range("D20").Value = $row21.D
$ws.Range("E20").Value = $row21.E
$ws.Range("F20").Value = $row21.F
$ws.Range("G20").Value = $row21.G
$ws.Range("H20").Value = $row21.H
$ws.Range("Q20").Value = $row21.Q
$ws.Range("R20").Value = $row21.R

$ws.Range("A21").Value = $row19.A
$ws.Range("B21").Value = $row19.B
$ws.Range("D21").Value = $row19.D
$ws.Range("E21").Value = $row19.E
$ws.Range("F21").Value = $row19.F
$ws.Range("G21").Value = $row19.G
$ws.Range("H21").Value = $row19.H
$ws.Range("Q21").Value = $row19.Q
$ws.Range("R21").Value = $row19.R
